$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read existing data (rows 2-23: Language, Value)
$data = @()
for ($r = 2; $r -le 23; $r++) {
    $lang = $ws.Cells.Item($r, 1).Value2
    $val = $ws.Cells.Item($r, 2).Value2
    $data += [PSCustomObject]@{Lang = $lang; Val = $val}
}

# Sort descending by value
$sorted = $data | Sort-Object -Property Val -Descending

# Drop the two lowest entries (Uzbek, Vietnamese)
$kept = $sorted[0..($sorted.Count - 3)]

# Write sorted data back into rows 2-21
$r = 2
foreach ($row in $kept) {
    $ws.Cells.Item($r, 1).Value = $row.Lang
    $ws.Cells.Item($r, 2).Value = $row.Val
    $r++
}

# Remove the now-obsolete rows 22 and 23
$ws.Range("A22:B23").Delete() | Out-Null
